$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '82.081.81'
$ws.Range("E2").Value = '  +2.84%  '

# Row 3
$ws.Range("D3").Value = '3.190.72'
$ws.Range("E3").Value = '  -0.42%  '

# Row 4
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").Value = '''215.97'
$ws.Range("E5").Value = '  +5.14%  '

# Row 6
$ws.Range("D6").Value = '''624.06'
$ws.Range("E6").Value = '  -1.56%  '

# Row 7
$ws.Range("E7").Value = '  +20.55%  '

# Row 8
$ws.Range("E8").Value = '  -0.08%  '

# Row 9
$ws.Range("D9").Value = '''0.586'
$ws.Range("E9").Value = '  +0.42%  '

# Row 10
$ws.Range("D10").Value = '3.185.10'
$ws.Range("E10").Value = '  -0.59%  '

# Row 11
$ws.Range("E11").Value = '  +1.88%  '

# Row 12
$ws.Range("E12").Value = '  +11.12%  '

# Row 13
$ws.Range("E13").Value = '  -0.29%  '

# Row 14
$ws.Range("E14").Value = '  -3.35%  '

# Row 15
$ws.Range("D15").Value = '3.775.34'
$ws.Range("E15").Value = '  -0.57%  '

# Row 16
$ws.Range("D16").Value = '''31.76'
$ws.Range("E16").Value = '  -0.73%  '

# Row 17
$ws.Range("D17").Value = '81.617.31'
$ws.Range("E17").Value = '  +2.49%  '

# Row 18
$ws.Range("D18").Value = '3.187.65'
$ws.Range("E18").Value = '  +0.06%  '

# Row 19
$ws.Range("D19").Value = '''3.22'
$ws.Range("E19").Value = '  +7.01%  '

# Row 20
$ws.Range("D20").Value = '''14.08'
$ws.Range("E20").Value = '  -3.05%  '

# Row 21
$ws.Range("D21").Value = '''435.87'
$ws.Range("E21").Value = '  +1.11%  '

# Row 22
$ws.Range("D22").Value = '''8.99'
$ws.Range("E22").Value = '  -2.22%  '

# Row 23
$ws.Range("D23").Value = '''5.15'
$ws.Range("E23").Value = '  +0.65%  '

# Row 24
$ws.Range("E24").Value = '  +5.92%  '

# Row 25
$ws.Range("D25").Value = '''5.35'
$ws.Range("E25").Value = '  +12.75%  '

# Row 26
$ws.Range("D26").Value = '3.358.04'
$ws.Range("E26").Value = '  -0.33%  '

# Row 27
$ws.Range("D27").Value = '''76.69'
$ws.Range("E27").Value = '  -0.79%  '

# Row 28
$ws.Range("D28").Value = '''11.05'
$ws.Range("E28").Value = '  -1.99%  '

# Row 29
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -0.13%  '

# Row 30
$ws.Range("E30").Value = '  +3.70%  '

# Row 31
$ws.Range("D31").Value = '''587.72'
$ws.Range("E31").Value = '  +11.42%  '

# Row 32
$ws.Range("D32").Value = '''9.10'
$ws.Range("E32").Value = '  +0.71%  '

# Row 33
$ws.Range("D33").Value = '''0.997'
$ws.Range("E33").Value = '  -0.22%  '

# Row 34
$ws.Range("E34").Value = '  +1.54%  '

# Row 35
$ws.Range("D35").Value = '''0.155'
$ws.Range("E35").Value = '  +8.94%  '

# Row 36
$ws.Range("D36").Value = '''2.01'
$ws.Range("E36").Value = '  +0.68%  '

# Row 37
$ws.Range("E37").Value = '  +17.14%  '

# Row 38
$ws.Range("E38").Value = '  -0.68%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").Value = '''6.17'
$ws.Range("E39").Value = '  +11.20%  '

# Row 40
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '''0.998'
$ws.Range("E40").Value = '  -0.14%  '

# Row 41
$ws.Range("E41").Value = '  +0.61%  '

# Row 42
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''2.05'
$ws.Range("E42").Value = '  +14.08%  '

# Row 43
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '''3.07'
$ws.Range("E43").Value = '  +22.21%  '

# Row 44
$ws.Range("E44").Value = '  +3.82%  '

# Row 45
$ws.Range("D45").Value = '''161.22'
$ws.Range("E45").Value = '  -2.52%  '

# Row 46
$ws.Range("E46").Value = '  +0.04%  '

# Row 47
$ws.Range("D47").Value = '''187.97'
$ws.Range("E47").Value = '  -2.55%  '

# Row 48
$ws.Range("D48").Value = '''44.74'
$ws.Range("E48").Value = '  +3.26%  '

# Row 49
$ws.Range("E49").Value = '  +1.01%  '

# Row 50
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''26.43'
$ws.Range("E50").Value = '  +2.05%  '

# Row 51
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '''0.775'
$ws.Range("E51").Value = '  -5.91%  '

